# This script applies updated NATMI TPM-derived metrics (commit: "update scripts wuth new tpm")
# to the Adam9-Itga3 ligand-receptor pair worksheet. Columns G:T (rows 2-17) are recomputed
# cell-by-cell for the new TPM inputs; columns A:F are identifiers/counts and are unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.573375
$ws.Range("H2").Value = 28.720125
$ws.Range("I2").Value = 0.1037691388643484
$ws.Range("J2").Value = 0.1037691388643484
$ws.Range("M2").Value = 11.25749966666667
$ws.Range("N2").Value = 33.772499
$ws.Range("O2").Value = 0.6929800609896341
$ws.Range("P2").Value = 0.6929800609896341
$ws.Range("Q2").Value = 107.772265871375
$ws.Range("R2").Value = 969.9503928423749
$ws.Range("S2").Value = 0.07190994417905795
$ws.Range("T2").Value = 0.07190994417905795

# Row 3
$ws.Range("G3").Value = 9.573375
$ws.Range("H3").Value = 28.720125
$ws.Range("I3").Value = 0.1037691388643484
$ws.Range("J3").Value = 0.1037691388643484
$ws.Range("M3").Value = 0.9898276666666668
$ws.Range("O3").Value = 0.06093101107050686
$ws.Range("P3").Value = 0.06093101107050686
$ws.Range("Q3").Value = 9.475991438375001
$ws.Range("R3").Value = 85.28392294537501
$ws.Range("S3").Value = 0.006322758548920575
$ws.Range("T3").Value = 0.006322758548920576

# Row 4
$ws.Range("G4").Value = 9.573375
$ws.Range("H4").Value = 28.720125
$ws.Range("I4").Value = 0.1037691388643484
$ws.Range("J4").Value = 0.1037691388643484
$ws.Range("M4").Value = 3.821582
$ws.Range("N4").Value = 11.464746
$ws.Range("O4").Value = 0.2352458543950409
$ws.Range("P4").Value = 0.2352458543950409
$ws.Range("Q4").Value = 36.58543757925
$ws.Range("R4").Value = 329.26893821325
$ws.Range("S4").Value = 0.02441125973198128
$ws.Range("T4").Value = 0.02441125973198128

# Row 5
$ws.Range("G5").Value = 9.573375
$ws.Range("H5").Value = 28.720125
$ws.Range("I5").Value = 0.1037691388643484
$ws.Range("J5").Value = 0.1037691388643484
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1761463333333333
$ws.Range("N5").Value = 0.528439
$ws.Range("O5").Value = 0.01084307354481826
$ws.Range("P5").Value = 0.01084307354481827
$ws.Range("Q5").Value = 1.686314903875
$ws.Range("R5").Value = 15.176834134875
$ws.Range("S5").Value = 0.001125176404388589
$ws.Range("T5").Value = 0.001125176404388589

# Row 6
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("I6").Value = 0.5196887643218222
$ws.Range("J6").Value = 0.5196887643218222
$ws.Range("M6").Value = 11.25749966666667
$ws.Range("N6").Value = 33.772499
$ws.Range("O6").Value = 0.6929800609896341
$ws.Range("P6").Value = 0.6929800609896341
$ws.Range("Q6").Value = 539.7369226709488
$ws.Range("R6").Value = 4857.632304038539
$ws.Range("S6").Value = 0.3601339515953639
$ws.Range("T6").Value = 0.3601339515953639

# Row 7
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("I7").Value = 0.5196887643218222
$ws.Range("J7").Value = 0.5196887643218222
$ws.Range("M7").Value = 0.9898276666666668
$ws.Range("O7").Value = 0.06093101107050686
$ws.Range("P7").Value = 0.06093101107050686
$ws.Range("Q7").Value = 47.45694466801812
$ws.Range("R7").Value = 427.1125020121631
$ws.Range("S7").Value = 0.03166516185211098
$ws.Range("T7").Value = 0.03166516185211098

# Row 8
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("I8").Value = 0.5196887643218222
$ws.Range("J8").Value = 0.5196887643218222
$ws.Range("M8").Value = 3.821582
$ws.Range("N8").Value = 11.464746
$ws.Range("O8").Value = 0.2352458543950409
$ws.Range("P8").Value = 0.2352458543950409
$ws.Range("Q8").Value = 183.2244254487674
$ws.Range("R8").Value = 1649.019829038906
$ws.Range("S8").Value = 0.1222546273823901
$ws.Range("T8").Value = 0.1222546273823901

# Row 9
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("I9").Value = 0.5196887643218222
$ws.Range("J9").Value = 0.5196887643218222
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1761463333333333
$ws.Range("N9").Value = 0.528439
$ws.Range("O9").Value = 0.01084307354481826
$ws.Range("P9").Value = 0.01084307354481827
$ws.Range("Q9").Value = 8.44527494631989
$ws.Range("R9").Value = 76.00747451687901
$ws.Range("S9").Value = 0.005635023491957243
$ws.Range("T9").Value = 0.005635023491957244

# Row 10
$ws.Range("G10").Value = 11.32006633333333
$ws.Range("H10").Value = 33.960199
$ws.Range("I10").Value = 0.122702133291269
$ws.Range("J10").Value = 0.122702133291269
$ws.Range("M10").Value = 11.25749966666667
$ws.Range("N10").Value = 33.772499
$ws.Range("O10").Value = 0.6929800609896341
$ws.Range("P10").Value = 0.6929800609896341
$ws.Range("Q10").Value = 127.4356429741446
$ws.Range("R10").Value = 1146.920786767301
$ws.Range("S10").Value = 0.08503013181174177
$ws.Range("T10").Value = 0.08503013181174177

# Row 11
$ws.Range("G11").Value = 11.32006633333333
$ws.Range("H11").Value = 33.960199
$ws.Range("I11").Value = 0.122702133291269
$ws.Range("J11").Value = 0.122702133291269
$ws.Range("M11").Value = 0.9898276666666668
$ws.Range("O11").Value = 0.06093101107050686
$ws.Range("P11").Value = 0.06093101107050686
$ws.Range("Q11").Value = 11.20491484523522
$ws.Range("R11").Value = 100.844233607117
$ws.Range("S11").Value = 0.007476365041945117
$ws.Range("T11").Value = 0.007476365041945118

# Row 12
$ws.Range("G12").Value = 11.32006633333333
$ws.Range("H12").Value = 33.960199
$ws.Range("I12").Value = 0.122702133291269
$ws.Range("J12").Value = 0.122702133291269
$ws.Range("M12").Value = 3.821582
$ws.Range("N12").Value = 11.464746
$ws.Range("O12").Value = 0.2352458543950409
$ws.Range("P12").Value = 0.2352458543950409
$ws.Range("Q12").Value = 43.26056173827267
$ws.Range("R12").Value = 389.345055644454
$ws.Range("S12").Value = 0.02886516818219875
$ws.Range("T12").Value = 0.02886516818219876

# Row 13
$ws.Range("G13").Value = 11.32006633333333
$ws.Range("H13").Value = 33.960199
$ws.Range("I13").Value = 0.122702133291269
$ws.Range("J13").Value = 0.122702133291269
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1761463333333333
$ws.Range("N13").Value = 0.528439
$ws.Range("O13").Value = 0.01084307354481826
$ws.Range("P13").Value = 0.01084307354481827
$ws.Range("Q13").Value = 1.993988177706778
$ws.Range("R13").Value = 17.945893599361
$ws.Range("S13").Value = 0.001330468255383323
$ws.Range("T13").Value = 0.001330468255383323

# Row 14
$ws.Range("G14").Value = 23.41838033333333
$ws.Range("H14").Value = 70.25514099999999
$ws.Range("I14").Value = 0.2538399635225604
$ws.Range("J14").Value = 0.2538399635225604
$ws.Range("M14").Value = 11.25749966666667
$ws.Range("N14").Value = 33.772499
$ws.Range("O14").Value = 0.6929800609896341
$ws.Range("P14").Value = 0.6929800609896341
$ws.Range("Q14").Value = 263.6324087963732
$ws.Range("R14").Value = 2372.691679167358
$ws.Range("S14").Value = 0.1759060334034704
$ws.Range("T14").Value = 0.1759060334034704

# Row 15
$ws.Range("G15").Value = 23.41838033333333
$ws.Range("H15").Value = 70.25514099999999
$ws.Range("I15").Value = 0.2538399635225604
$ws.Range("J15").Value = 0.2538399635225604
$ws.Range("M15").Value = 0.9898276666666668
$ws.Range("O15").Value = 0.06093101107050686
$ws.Range("P15").Value = 0.06093101107050686
$ws.Range("Q15").Value = 23.18016076245589
$ws.Range("R15").Value = 208.621446862103
$ws.Range("S15").Value = 0.01546672562753019
$ws.Range("T15").Value = 0.01546672562753019

# Row 16
$ws.Range("G16").Value = 23.41838033333333
$ws.Range("H16").Value = 70.25514099999999
$ws.Range("I16").Value = 0.2538399635225604
$ws.Range("J16").Value = 0.2538399635225604
$ws.Range("M16").Value = 3.821582
$ws.Range("N16").Value = 11.464746
$ws.Range("O16").Value = 0.2352458543950409
$ws.Range("P16").Value = 0.2352458543950409
$ws.Range("Q16").Value = 89.49526075102065
$ws.Range("R16").Value = 805.457346759186
$ws.Range("S16").Value = 0.05971479909847074
$ws.Range("T16").Value = 0.05971479909847074

# Row 17
$ws.Range("G17").Value = 23.41838033333333
$ws.Range("H17").Value = 70.25514099999999
$ws.Range("I17").Value = 0.2538399635225604
$ws.Range("J17").Value = 0.2538399635225604
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1761463333333333
$ws.Range("N17").Value = 0.528439
$ws.Range("O17").Value = 0.01084307354481826
$ws.Range("P17").Value = 0.01084307354481827
$ws.Range("Q17").Value = 4.12506182832211
$ws.Range("R17").Value = 37.125556454899
$ws.Range("S17").Value = 0.002752405393089108
$ws.Range("T17").Value = 0.002752405393089108

Write-Output "Applied 192 cell updates to Sheet1!G2:T17"
